{"js": "// Update the statistics table in the \"Open Data Suchportal\" overview:\n// geoportal/portal counts increased, so bump the corresponding figures.\nconst replacements = [\n  // [searchText, replacementText]\n  [\"(90)\", \"(91)\"],\n  [\" (64)\", \" (65)\"],\n  [\"      10\", \"      11\"],\n  [\"(74)\", \"(77)\"],\n  [\" (59)\", \" (62)\"],\n  [\"      33\", \"      32\"],\n  [\"      24\", \"      28\"],\n  [\"Datenportale insgesamt: 164\", \"Datenportale insgesamt: 168\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${searchText}\"`);\n  }\n\n  // Each search term is unique in this document, so replace the single match.\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the statistics table in the \"Open Data Suchportal\" overview:\n# geoportal/portal counts increased, so bump the corresponding figures.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"(90)\";  Replace = \"(91)\" },\n    @{ Find = \" (64)\"; Replace = \" (65)\" },\n    @{ Find = \"      10\"; Replace = \"      11\" },\n    @{ Find = \"(74)\";  Replace = \"(77)\" },\n    @{ Find = \" (59)\"; Replace = \" (62)\" },\n    @{ Find = \"      33\"; Replace = \"      32\" },\n    @{ Find = \"      24\"; Replace = \"      28\" },\n    @{ Find = \"Datenportale insgesamt: 164\"; Replace = \"Datenportale insgesamt: 168\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
